# daily auto push: 2026-02-10 19:22 UTC
# Insert one new data row (2026/02/11) right after the existing last
# 2026/02/10 entry (row 806), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 806..847 down to 807..848, leaving a blank row 806.
$ws.Rows.Item(806).EntireRow.Insert()

# Force column A to Text format first so the "yyyy/mm/dd"-looking string
# is stored verbatim instead of being auto-parsed into a date serial.
$ws.Range("A806").NumberFormat = "@"
$ws.Range("A806").Value = "2026/02/11"
$ws.Range("B806").Value = "水"
$ws.Range("C806").Value = 0
$ws.Range("D806").Value = 201
